$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Phase3")

# Update column E (rows 2-18) from 0.6 to 0.75. Column F holds =LN(E#) formulas
# that will recalculate automatically.
$ws.Range("E2:E18").Value = 0.75

# Update the selected cell on this sheet to E24 (matches final saved selection).
$ws.Activate()
$ws.Range("E24").Select()
